{"js": "// Remove the two empty paragraphs that immediately follow the paragraph\n// \"Cualquier decisi\u00f3n del proyecto se tiene que tomar en conjunto.\" \u2014\n// one empty numbered list item and one empty plain paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Cualquier decisi\u00f3n del proyecto se tiene que tomar en conjunto.\";\nlet markerIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex === -1) {\n  throw new Error(\"Could not locate the target paragraph: \" + marker);\n}\n\nconst firstEmpty = paragraphs.items[markerIndex + 1];\nconst secondEmpty = paragraphs.items[markerIndex + 2];\n\nfirstEmpty.delete();\nsecondEmpty.delete();\n\nawait context.sync();\n", "ps1": "# Remove the two empty paragraphs that immediately follow the paragraph\n# \"Cualquier decisi\u00f3n del proyecto se tiene que tomar en conjunto.\" \u2014\n# one empty numbered list item and one empty plain paragraph.\n$d = $word.ActiveDocument\n\n$marker = \"Cualquier decisi\u00f3n del proyecto se tiene que tomar en conjunto.\"\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $marker) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the target paragraph: $marker\"\n}\n\n# Both paragraphs being removed land at the same index (targetIndex + 1)\n# once the prior one is deleted, since later paragraphs shift up.\n$d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n$d.Paragraphs.Item($targetIndex + 1).Range.Delete()\n"}
